# Replace the four "statut" marker values (column A) with their new
# equivalents, per the emoji -> plain-sign / warning-sign rework:
#   book(red)    "📕" -> "-3"
#   book(blue)   "📘" -> "⚠️"
#   book(orange) "📙" -> "+3"
#   book(green)  "📗" -> "✅"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$lastRow = $used.Rows.Count + $used.Row - 1

$map = @{
    "📕" = "-3"
    "📘" = "⚠️"
    "📙" = "+3"
    "📗" = "✅"
}
# Replacement values that Excel would otherwise auto-convert to numbers
# need to be forced to text so they round-trip as strings, not numerics.
$numericLooking = @{ "-3" = $true; "+3" = $true }

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $old = $cell.Value2
    if ($map.ContainsKey($old)) {
        $new = $map[$old]
        if ($numericLooking.ContainsKey($new)) {
            $cell.NumberFormat = "@"
        }
        $cell.Value2 = $new
    }
}
